$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.115.26'
$ws.Range("E2").Value = '  +1.48%  '

$ws.Range("D3").Value = '1.931.00'
$ws.Range("E3").Value = '  +2.38%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.66'
$ws.Range("E5").Value = '  +1.49%  '

$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4614'
$ws.Range("E7").Value = '  +1.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3834'
$ws.Range("E8").Value = '  +1.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07777'
$ws.Range("E9").Value = '  +1.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9808'
$ws.Range("E10").Value = '  +2.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.60'
$ws.Range("E11").Value = '  +3.36%  '

$ws.Range("D12").Value = '1.960.77'
$ws.Range("E12").Value = '  +3.90%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.708'
$ws.Range("E13").Value = '  +1.40%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.985'
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07053'
$ws.Range("E15").Value = '  +0.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.49'
$ws.Range("E17").Value = '  +2.62%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009556'
$ws.Range("E18").Value = '  +1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.78'
$ws.Range("E19").Value = '  +1.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("D21").Value = '29.122.49'
$ws.Range("E21").Value = '  +1.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.355'
$ws.Range("E22").Value = '  +1.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.97'
$ws.Range("E23").Value = '  +1.42%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.092'
$ws.Range("E24").Value = '  +0.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.28'
$ws.Range("E25").Value = '  +2.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.09'
$ws.Range("E26").Value = '  +1.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.684'
$ws.Range("E27").Value = '  +1.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '118.22'
$ws.Range("E28").Value = '  +1.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.852'
$ws.Range("E29").Value = '  +2.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09349'
$ws.Range("E30").Value = '  +1.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8671'
$ws.Range("E31").Value = '  +3.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.151'
$ws.Range("E32").Value = '  +2.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.253'
$ws.Range("E33").Value = '  +1.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.018'
$ws.Range("E34").Value = '  -1.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05708'
$ws.Range("E35").Value = '  +1.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.161'

$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02057'
$ws.Range("E38").Value = '  +1.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.084'
$ws.Range("E39").Value = '  +14.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.563'
$ws.Range("E40").Value = '  +1.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5529'
$ws.Range("E41").Value = '  +1.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1759'
$ws.Range("E42").Value = '  +1.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.380'
$ws.Range("E43").Value = '  +2.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000002843'
$ws.Range("E44").Value = '  -2.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.206'
$ws.Range("E45").Value = '  +6.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5212'
$ws.Range("E46").Value = '  +1.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.28'
$ws.Range("E47").Value = '  +0.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06932'
$ws.Range("E48").Value = '  +2.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.781'
$ws.Range("E49").Value = '  +1.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.45'
$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.32%  '
